# "chose US for sprint 3"
# Assign user stories US31, US35 and US39 to Sprint 3 in the Backlog sheet,
# and add the newly-chosen stories (US35, US39) with their sub-tasks to the
# Sprint3 sheet (US31 was already present on the Sprint3 sheet).

$wb = $excel.ActiveWorkbook

# ---- Backlog sheet: mark Sprint = 3 for US31, US35 and US39 ----
$backlog = $wb.Worksheets.Item("Backlog")

# US31 - "List living single" - just gets assigned to sprint 3
$backlog.Cells.Item(32, 1).Value = 3

# US35 - "List recent births" - assigned to sprint 3, owned by Maha, Initiated
$backlog.Cells.Item(36, 1).Value = 3
$backlog.Cells.Item(36, 4).Value = "MA"
$backlog.Cells.Item(36, 5).Value = "Initiated"

# US39 - "List upcoming anniversaries" - assigned to sprint 3, owned by Maha, Initiated
$backlog.Cells.Item(40, 1).Value = 3
$backlog.Cells.Item(40, 4).Value = "MA"
$backlog.Cells.Item(40, 5).Value = "Initiated"

# ---- Sprint3 sheet: add the newly chosen stories with their sub-tasks ----
$sprint3 = $wb.Worksheets.Item("Sprint3")

# US35 - List recent births
$sprint3.Cells.Item(14, 1).Value = "US35"
$sprint3.Cells.Item(14, 1).Font.Bold = $true
$sprint3.Cells.Item(14, 2).Value = "List recent births"
$sprint3.Cells.Item(14, 3).Value = "Maha"
$sprint3.Cells.Item(14, 4).Value = "Initiated"
$sprint3.Cells.Item(14, 5).Value = 30
$sprint3.Cells.Item(14, 6).Value = 60

$sprint3.Cells.Item(16, 2).Value = "find birthdate"
$sprint3.Cells.Item(17, 2).Value = "compare birth date to today"
$sprint3.Cells.Item(18, 2).Value = "if less than # months print Name"

$sprint3.Cells.Item(16, 1).Value = "T35.01"
$sprint3.Cells.Item(17, 1).Value = "T35.02"
$sprint3.Cells.Item(18, 1).Value = "T35.03"

# US39 - List upcoming anniversaries
$sprint3.Cells.Item(21, 1).Value = "US39"
$sprint3.Cells.Item(21, 1).Font.Bold = $true
$sprint3.Cells.Item(21, 2).Value = "List upcoming anniversaries"
$sprint3.Cells.Item(21, 3).Value = "Maha"
$sprint3.Cells.Item(21, 4).Value = "Initiated"
$sprint3.Cells.Item(21, 5).Value = 35
$sprint3.Cells.Item(21, 6).Value = 60

$sprint3.Cells.Item(23, 2).Value = "Find Marriage date for each family"
$sprint3.Cells.Item(24, 2).Value = "compare to today"
$sprint3.Cells.Item(25, 2).Value = "print marriage date for next month"

$sprint3.Cells.Item(23, 1).Value = "T39.01"
$sprint3.Cells.Item(24, 1).Value = "T39.02"
$sprint3.Cells.Item(25, 1).Value = "T39.03"
